$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the styled format (bold, border, centered) from an existing column-A data cell
# so newly appended rows (38, 39) match the existing column-A formatting.
$aStyle = $ws.Cells.Item(2, 1).Style
$ws.Cells.Item(38, 1).Style = $aStyle
$ws.Cells.Item(39, 1).Style = $aStyle

# Row 2
$ws.Cells.Item(2, 1).Value = "'14"
$ws.Cells.Item(2, 2).Value = 5809570196
$ws.Cells.Item(2, 3).Value = "2021-11-21 23:33:21"
$ws.Cells.Item(2, 4).Value = "桃树下的孩子"
$ws.Cells.Item(2, 5).Value = "感谢字幕！第六场的普通版剪辑版双声道的少年/少女记忆虽然别有一番风味，而且剪辑在一起能明显看出很多动作都是同步或者对称的，对照着看特别有感觉，但是能看到清晰独立版本的活着的只有我（？）真的好棒！以及！！！真诚安利大家关注这场犬彦和宫比的互动，包括狼欒神社solo的时候模拟打鼓啦，井户曲摸头啦，魔神曲犬彦拉住宫比的手揽住他的腰等等，互动又多又甜（当然其他场也很甜！我记得有一场，忘记是哪场了狼欒solo时宫比在后面给犬彦比心心！），这对青梅竹马（大概是吧）真的超级超级好嗑！！！"
$ws.Cells.Item(2, 6).Value = 5
$ws.Cells.Item(2, 7).Value = 1
$ws.Cells.Item(2, 8).Value = 0

# Row 3
$ws.Cells.Item(3, 1).Value = "'14.1"
$ws.Cells.Item(3, 2).Value = 5809780005
$ws.Cells.Item(3, 3).Value = "2021-11-22 00:05:05"
$ws.Cells.Item(3, 4).Value = "Ponster_"
$ws.Cells.Item(3, 5).Value = "弹幕中的翻译佬！！感谢指正[脱单doge]"
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 5809570196

# Row 4
$ws.Cells.Item(4, 1).Value = "'15"
$ws.Cells.Item(4, 2).Value = 5809565870
$ws.Cells.Item(4, 3).Value = "2021-11-21 23:33:05"
$ws.Cells.Item(4, 4).Value = "我心向云月"
$ws.Cells.Item(4, 5).Value = "西装男跟夜姬一起演太可怕了，jk快远离疯批男[冷]"
$ws.Cells.Item(4, 6).Value = 4
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0

# Row 5
$ws.Cells.Item(5, 1).Value = "'1"
$ws.Cells.Item(5, 2).Value = 5819891283
$ws.Cells.Item(5, 3).Value = "2021-11-23 20:01:10"
$ws.Cells.Item(5, 4).Value = "霜小凝"
$ws.Cells.Item(5, 5).Value = "一键三连了！"
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0

# Row 6
$ws.Cells.Item(6, 1).Value = "'12"
$ws.Cells.Item(6, 2).Value = 5809727182
$ws.Cells.Item(6, 3).Value = "2021-11-21 23:58:18"
$ws.Cells.Item(6, 4).Value = "夏空凛冬至"
$ws.Cells.Item(6, 5).Value = "谢谢up！不知道还有没有其他的！"
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 3
$ws.Cells.Item(6, 8).Value = 0

# Row 7
$ws.Cells.Item(7, 1).Value = "'12.1"
$ws.Cells.Item(7, 2).Value = 5809752664
$ws.Cells.Item(7, 3).Value = "2021-11-22 00:01:40"
$ws.Cells.Item(7, 4).Value = "Ponster_"
$ws.Cells.Item(7, 5).Value = "暂时...不会做了吧，原因看视频开头"
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 5809727182

# Row 8
$ws.Cells.Item(8, 1).Value = "'12.2"
$ws.Cells.Item(8, 2).Value = 5809822355
$ws.Cells.Item(8, 3).Value = "2021-11-22 00:13:14"
$ws.Cells.Item(8, 4).Value = "我心向云月"
$ws.Cells.Item(8, 5).Value = "回复 @Ponster_ :可以求其他安可的生肉吗[大哭]"
$ws.Cells.Item(8, 6).Value = 0
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 5809727182

# Row 9
$ws.Cells.Item(9, 1).Value = "'12.3"
$ws.Cells.Item(9, 2).Value = 5810008526
$ws.Cells.Item(9, 3).Value = "2021-11-22 00:49:30"
$ws.Cells.Item(9, 4).Value = "Hexachlorocyclohexane"
$ws.Cells.Item(9, 5).Value = "回复 @Ponster_ :同求其他安可的生肉[大哭]"
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 5809727182

# Row 10
$ws.Cells.Item(10, 1).Value = "'20"
$ws.Cells.Item(10, 2).Value = 5807519793
$ws.Cells.Item(10, 3).Value = "2021-11-21 18:25:36"
$ws.Cells.Item(10, 4).Value = "墨弦青风"
$ws.Cells.Item(10, 5).Value = "感谢up，up辛苦了[热词系列_吹爆]"
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 0

# Row 11
$ws.Cells.Item(11, 1).Value = "'20.1"
$ws.Cells.Item(11, 2).Value = 5810526423
$ws.Cells.Item(11, 3).Value = "2021-11-22 06:24:28"
$ws.Cells.Item(11, 4).Value = "Ponster_"
$ws.Cells.Item(11, 5).Value = "感谢(=・ω・=)"
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 5807519793

# Row 12
$ws.Cells.Item(12, 1).Value = "'19"
$ws.Cells.Item(12, 2).Value = 5808208456
$ws.Cells.Item(12, 3).Value = "2021-11-21 20:10:25"
$ws.Cells.Item(12, 4).Value = "Ponster_"
$ws.Cells.Item(12, 5).Value = "前面传错了版本，已更正。`n曲目信息、想说的话都在视频里了。`n第一场重制版也已上传，链接https://www.bilibili.com/video/BV1EU4y1u7HA"
$ws.Cells.Item(12, 6).Value = 0
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0

# Row 13
$ws.Cells.Item(13, 1).Value = "'18"
$ws.Cells.Item(13, 2).Value = 5808292895
$ws.Cells.Item(13, 3).Value = "2021-11-21 20:23:50"
$ws.Cells.Item(13, 4).Value = "烟云z"
$ws.Cells.Item(13, 5).Value = "太顶了老哥[BW2020_棒棒哦]"
$ws.Cells.Item(13, 6).Value = 0
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0

# Row 14
$ws.Cells.Item(14, 1).Value = "'17"
$ws.Cells.Item(14, 2).Value = 5808581488
$ws.Cells.Item(14, 3).Value = "2021-11-21 21:09:35"
$ws.Cells.Item(14, 4).Value = "羽蛇的尾巴尖"
$ws.Cells.Item(14, 5).Value = "迅速缓存"
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = 0

# Row 15
$ws.Cells.Item(15, 1).Value = "'17.1"
$ws.Cells.Item(15, 2).Value = 5810522401
$ws.Cells.Item(15, 3).Value = "2021-11-22 06:24:00"
$ws.Cells.Item(15, 4).Value = "Ponster_"
$ws.Cells.Item(15, 5).Value = "[tv_点赞]"
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 5808581488

# Row 16
$ws.Cells.Item(16, 1).Value = "'16"
$ws.Cells.Item(16, 2).Value = 5809188522
$ws.Cells.Item(16, 3).Value = "2021-11-21 22:40:13"
$ws.Cells.Item(16, 4).Value = "Panic-"
$ws.Cells.Item(16, 5).Value = "草西装男好可怕"
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0

# Row 17
$ws.Cells.Item(17, 1).Value = "'13"
$ws.Cells.Item(17, 2).Value = 5809720823
$ws.Cells.Item(17, 3).Value = "2021-11-21 23:56:28"
$ws.Cells.Item(17, 4).Value = "VirginMary"
$ws.Cells.Item(17, 5).Value = "好耶 来了来了"
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0

# Row 18
$ws.Cells.Item(18, 1).Value = "'11"
$ws.Cells.Item(18, 2).Value = 5810353355
$ws.Cells.Item(18, 3).Value = "2021-11-22 03:00:50"
$ws.Cells.Item(18, 4).Value = "召唤魔术"
$ws.Cells.Item(18, 5).Value = "求上传其他场次生肉[tv_大佬]"
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0

# Row 19
$ws.Cells.Item(19, 1).Value = "'2"
$ws.Cells.Item(19, 2).Value = 5817605996
$ws.Cells.Item(19, 3).Value = "2021-11-23 12:32:48"
$ws.Cells.Item(19, 4).Value = "总攻祁墨宸大人"
$ws.Cells.Item(19, 5).Value = "飞速闻讯而来"
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0

# Row 20
$ws.Cells.Item(20, 1).Value = "'10"
$ws.Cells.Item(20, 2).Value = 5810538568
$ws.Cells.Item(20, 3).Value = "2021-11-22 06:45:54"
$ws.Cells.Item(20, 4).Value = "冬の伝言"
$ws.Cells.Item(20, 5).Value = "唉，怎么说呢，自从进击的轨迹之后近年的陛下仿佛是换了一种形象，虽说是放开了许多，但也让人感到缺失了5.6.7平表演时的悲伤、感动和纪行时的那种温柔、坚毅、富有人格魅力的感觉"
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 3
$ws.Cells.Item(20, 8).Value = 0

# Row 21
$ws.Cells.Item(21, 1).Value = "'10.1"
$ws.Cells.Item(21, 2).Value = 5811776235
$ws.Cells.Item(21, 3).Value = "2021-11-22 12:37:31"
$ws.Cells.Item(21, 4).Value = "不是你的朱雀"
$ws.Cells.Item(21, 5).Value = "从忧郁小王子变成了阳光大男孩"
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 5810538568

# Row 22
$ws.Cells.Item(22, 1).Value = "'10.2"
$ws.Cells.Item(22, 2).Value = 5814346853
$ws.Cells.Item(22, 3).Value = "2021-11-22 20:52:20"
$ws.Cells.Item(22, 4).Value = "Ponster_"
$ws.Cells.Item(22, 5).Value = "说起来陛下已经是40+的中年大叔了啊（大不敬）"
$ws.Cells.Item(22, 6).Value = 0
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 5810538568

# Row 23
$ws.Cells.Item(23, 1).Value = "'10.3"
$ws.Cells.Item(23, 2).Value = 5818583189
$ws.Cells.Item(23, 3).Value = "2021-11-23 16:14:00"
$ws.Cells.Item(23, 4).Value = "燈留子"
$ws.Cells.Item(23, 5).Value = "陛下2040+岁了变得开朗也很不错"
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 5810538568

# Row 24
$ws.Cells.Item(24, 1).Value = "'9"
$ws.Cells.Item(24, 2).Value = 5812089382
$ws.Cells.Item(24, 3).Value = "2021-11-22 13:36:01"
$ws.Cells.Item(24, 4).Value = "enemin"
$ws.Cells.Item(24, 5).Value = "太感谢了 之前还有一版安可不知道有没有大佬传 几乎是猫咪铃唱人偶的"
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 3
$ws.Cells.Item(24, 8).Value = 0

# Row 25
$ws.Cells.Item(25, 1).Value = "'9.1"
$ws.Cells.Item(25, 2).Value = 5814158584
$ws.Cells.Item(25, 3).Value = "2021-11-22 20:21:43"
$ws.Cells.Item(25, 4).Value = "林花花花"
$ws.Cells.Item(25, 5).Value = "那个b站之前有，是被删了吗"
$ws.Cells.Item(25, 6).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 5812089382

# Row 26
$ws.Cells.Item(26, 1).Value = "'9.2"
$ws.Cells.Item(26, 2).Value = 5814154479
$ws.Cells.Item(26, 3).Value = "2021-11-22 20:22:06"
$ws.Cells.Item(26, 4).Value = "林花花花"
$ws.Cells.Item(26, 5).Value = "还好缓存的快"
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 8).Value = 5812089382

# Row 27
$ws.Cells.Item(27, 1).Value = "'9.3"
$ws.Cells.Item(27, 2).Value = 5814448893
$ws.Cells.Item(27, 3).Value = "2021-11-22 21:08:40"
$ws.Cells.Item(27, 4).Value = "enemin"
$ws.Cells.Item(27, 5).Value = "回复 @林花花花 :是啊 我前一秒还在看 然后推出去发现就无了 没有缓存 伤心了"
$ws.Cells.Item(27, 6).Value = 0
$ws.Cells.Item(27, 7).Value = 0
$ws.Cells.Item(27, 8).Value = 5812089382

# Row 28
$ws.Cells.Item(28, 1).Value = "'8"
$ws.Cells.Item(28, 2).Value = 5812890404
$ws.Cells.Item(28, 3).Value = "2021-11-22 16:48:22"
$ws.Cells.Item(28, 4).Value = "大白梨°"
$ws.Cells.Item(28, 5).Value = "[doge]陛下终究还是老了，不装13我很不习惯的"
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 0

# Row 29
$ws.Cells.Item(29, 1).Value = "'7"
$ws.Cells.Item(29, 2).Value = 5813229158
$ws.Cells.Item(29, 3).Value = "2021-11-22 17:57:18"
$ws.Cells.Item(29, 4).Value = "林花花花"
$ws.Cells.Item(29, 5).Value = "感谢！！"
$ws.Cells.Item(29, 6).Value = 0
$ws.Cells.Item(29, 7).Value = 2
$ws.Cells.Item(29, 8).Value = 0

# Row 30
$ws.Cells.Item(30, 1).Value = "'7.1"
$ws.Cells.Item(30, 2).Value = 5813402664
$ws.Cells.Item(30, 3).Value = "2021-11-22 18:25:22"
$ws.Cells.Item(30, 4).Value = "Ponster_"
$ws.Cells.Item(30, 5).Value = "也感谢你(=・ω・=)"
$ws.Cells.Item(30, 6).Value = 0
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 5813229158

# Row 31
$ws.Cells.Item(31, 1).Value = "'7.2"
$ws.Cells.Item(31, 2).Value = 5814183617
$ws.Cells.Item(31, 3).Value = "2021-11-22 20:26:08"
$ws.Cells.Item(31, 4).Value = "林花花花"
$ws.Cells.Item(31, 5).Value = "回复 @Ponster_ :想看八平[笑哭]"
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 5813229158

# Row 32
$ws.Cells.Item(32, 1).Value = "'6"
$ws.Cells.Item(32, 2).Value = 5815080989
$ws.Cells.Item(32, 3).Value = "2021-11-22 22:43:19"
$ws.Cells.Item(32, 4).Value = "艾奥萝卜"
$ws.Cells.Item(32, 5).Value = "先马再看"
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 0

# Row 33
$ws.Cells.Item(33, 1).Value = "'5"
$ws.Cells.Item(33, 2).Value = 5815188052
$ws.Cells.Item(33, 3).Value = "2021-11-22 22:59:05"
$ws.Cells.Item(33, 4).Value = "自律-Official"
$ws.Cells.Item(33, 5).Value = "哦哦哦赶紧缓存爽到[夏诺雅_太会了]"
$ws.Cells.Item(33, 6).Value = 0
$ws.Cells.Item(33, 7).Value = 2
$ws.Cells.Item(33, 8).Value = 0

# Row 34
$ws.Cells.Item(34, 1).Value = "'5.1"
$ws.Cells.Item(34, 2).Value = 5815284440
$ws.Cells.Item(34, 3).Value = "2021-11-22 23:12:58"
$ws.Cells.Item(34, 4).Value = "Ponster_"
$ws.Cells.Item(34, 5).Value = "这里不要脸地对简介链接里的仓库进行一个安利[doge]`n视频的评论区、弹幕已备份到其中，刚刚更新过`n[吃瓜][吃瓜][吃瓜]"
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 0
$ws.Cells.Item(34, 8).Value = 5815188052

# Row 35
$ws.Cells.Item(35, 1).Value = "'5.2"
$ws.Cells.Item(35, 2).Value = 5815389687
$ws.Cells.Item(35, 3).Value = "2021-11-22 23:28:15"
$ws.Cells.Item(35, 4).Value = "自律-Official"
$ws.Cells.Item(35, 5).Value = "回复 @Ponster_ :好起来了[夏诺雅_震撼]"
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 0
$ws.Cells.Item(35, 8).Value = 5815188052

# Row 36
$ws.Cells.Item(36, 1).Value = "'4"
$ws.Cells.Item(36, 2).Value = 5815263944
$ws.Cells.Item(36, 3).Value = "2021-11-22 23:09:47"
$ws.Cells.Item(36, 4).Value = "敦肃皇贵妃葛小队"
$ws.Cells.Item(36, 5).Value = "第一次看西装男吓尿了"
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 1
$ws.Cells.Item(36, 8).Value = 0

# Row 37
$ws.Cells.Item(37, 1).Value = "'4.1"
$ws.Cells.Item(37, 2).Value = 5815293406
$ws.Cells.Item(37, 3).Value = "2021-11-22 23:14:02"
$ws.Cells.Item(37, 4).Value = "Ponster_"
$ws.Cells.Item(37, 5).Value = "陛下可是“音乐界的杀人贵公子”呢"
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 5815263944

# Row 38
$ws.Cells.Item(38, 1).Value = "'3"
$ws.Cells.Item(38, 2).Value = 5815985159
$ws.Cells.Item(38, 3).Value = "2021-11-23 01:26:12"
$ws.Cells.Item(38, 4).Value = "木容秀吉"
$ws.Cells.Item(38, 5).Value = "这个安可我直呼好家伙"
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 0

# Row 39
$ws.Cells.Item(39, 1).Value = "'21"
$ws.Cells.Item(39, 2).Value = 5807507484
$ws.Cells.Item(39, 3).Value = "2021-11-21 18:22:28"
$ws.Cells.Item(39, 4).Value = "syyuansang"
$ws.Cells.Item(39, 5).Value = "这是第几场，泪目了"
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(39, 7).Value = 0
$ws.Cells.Item(39, 8).Value = 0
